$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns at D:F, shifting the old "Terms Typically Offered"
# column (D) to G.
$ws.Range("D1:F1").EntireColumn.Insert()

# New header row values for the inserted columns.
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

# Fill "NA" for each data row (2-13) in the three new columns.
$ws.Range("D2:F13").Value = "NA"
